# ---------------------------------------------------------------------------
# edit.ps1
#
# Reproduces the "Add files via upload" commit: extends the YTY ("Yatay
# Tambur Yaylandirma" / horizontal-drum straightening line) machine schema
# in the fields/options lookup sheets, and leaves the workbook's view state
# (active sheet + selections) the way the author left it when they saved.
#
# New field groups added to the "fields" sheet (rows 28-35):
#   - Discap_YTY            (roll outer diameter select + "other" text field)
#   - Dogrultma_Tipi_YTY    (straightener type, gated on Makina_Tipi = CMDC)
#   - Giris_Unitesi_YTY     (entry unit, gated on Makina_Tipi = CMDC)
#   - Hiz_YTY               (opener/straightener speed selects + "other" text)
#
# New option lists added to the "options" sheet (rows 86-98) backing the
# *_opts OptionsKeys referenced above.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# Auto-generated cell writes for sheet3 ("fields") rows 28-35
$ws3 = $wb.Worksheets.Item("fields")
$ws3.Cells.Item(28, 1).Value = "Discap_YTY"
$ws3.Cells.Item(28, 2).Value = "Discap_YTY"
$ws3.Cells.Item(28, 3).Value = "Rulo Dış Çapı"
$ws3.Cells.Item(28, 4).Value = "select"
$ws3.Cells.Item(28, 5).Value = $true
$ws3.Cells.Item(28, 6).Value = "Discap_YTY_opts"
$ws3.Cells.Item(28, 17).Value = "radio"
$ws3.Cells.Item(29, 1).Value = "Discap_YTY"
$ws3.Cells.Item(29, 2).Value = "Discap_harici"
$ws3.Cells.Item(29, 3).Value = "Listede Olmayan Rulo Dış Çapı (belirtiniz...)"
$ws3.Cells.Item(29, 4).Value = "number"
$ws3.Cells.Item(29, 5).Value = $false
$ws3.Cells.Item(29, 6).Value = "Discap_harici_opts"
$ws3.Cells.Item(29, 11).Value = "D"
$ws3.Cells.Item(29, 13).Value = 4
$ws3.Cells.Item(29, 14).Value = "Discap_YTY"
$ws3.Cells.Item(29, 15).Value = "Diğer"
$ws3.Cells.Item(29, 16).Value = 0
$ws3.Cells.Item(30, 1).Value = "Dogrultma_Tipi_YTY"
$ws3.Cells.Item(30, 2).Value = "Dogrultma_Tipi_YTY"
$ws3.Cells.Item(30, 3).Value = "Doğrultma Tipi"
$ws3.Cells.Item(30, 4).Value = "select"
$ws3.Cells.Item(30, 5).Value = $false
$ws3.Cells.Item(30, 6).Value = "Dogrultma_Tipi_YTY_opts"
$ws3.Cells.Item(30, 14).Value = "Makina_Tipi"
$ws3.Cells.Item(30, 15).Value = "CMDC"
$ws3.Cells.Item(30, 17).Value = "radio"
$ws3.Cells.Item(31, 1).Value = "Giris_Unitesi_YTY"
$ws3.Cells.Item(31, 2).Value = "Giris_Unitesi_YTY"
$ws3.Cells.Item(31, 3).Value = "Giriş Ünitesi"
$ws3.Cells.Item(31, 4).Value = "select"
$ws3.Cells.Item(31, 5).Value = $false
$ws3.Cells.Item(31, 6).Value = "Giris_Unitesi_YTY_opts"
$ws3.Cells.Item(31, 14).Value = "Makina_Tipi"
$ws3.Cells.Item(31, 15).Value = "CMDC"
$ws3.Cells.Item(31, 17).Value = "radio"
$ws3.Cells.Item(32, 1).Value = "Hiz_YTY"
$ws3.Cells.Item(32, 2).Value = "Hiz_acici_YTY"
$ws3.Cells.Item(32, 3).Value = "Tambur Hızı"
$ws3.Cells.Item(32, 4).Value = "select"
$ws3.Cells.Item(32, 5).Value = $true
$ws3.Cells.Item(32, 6).Value = "Hiz_acici_YTY_opts"
$ws3.Cells.Item(32, 17).Value = "radio"
$ws3.Cells.Item(33, 1).Value = "Hiz_YTY"
$ws3.Cells.Item(33, 2).Value = "Hiz_acici_harici"
$ws3.Cells.Item(33, 3).Value = "Listede Olmayan Açıcı Hızı (belirtiniz...)"
$ws3.Cells.Item(33, 4).Value = "text"
$ws3.Cells.Item(33, 5).Value = $false
$ws3.Cells.Item(33, 6).Value = "Hiz_acici_harici_opts"
$ws3.Cells.Item(33, 14).Value = "Hiz_acici_YTY"
$ws3.Cells.Item(33, 15).Value = "Diğer"
$ws3.Cells.Item(34, 1).Value = "Hiz_YTY"
$ws3.Cells.Item(34, 2).Value = "Hiz_dogrultucu_YTY"
$ws3.Cells.Item(34, 3).Value = "Doğrultucu Hızı"
$ws3.Cells.Item(34, 4).Value = "select"
$ws3.Cells.Item(34, 5).Value = $false
$ws3.Cells.Item(34, 6).Value = "Hiz_dogrultucu_YTY_opts"
$ws3.Cells.Item(34, 14).Value = "Makina_Tipi"
$ws3.Cells.Item(34, 15).Value = "CMDC"
$ws3.Cells.Item(34, 17).Value = "radio"
$ws3.Cells.Item(35, 1).Value = "Hiz_YTY"
$ws3.Cells.Item(35, 2).Value = "Hiz_dogrultucu_harici"
$ws3.Cells.Item(35, 3).Value = "Listede Olmayan Doğrultma Hızı (belirtiniz...)"
$ws3.Cells.Item(35, 4).Value = "text"
$ws3.Cells.Item(35, 5).Value = $false
$ws3.Cells.Item(35, 6).Value = "Hiz_dogrultucu_harici_opts"
$ws3.Cells.Item(35, 14).Value = "Hiz_dogrultucu_YTY"
$ws3.Cells.Item(35, 15).Value = "Diğer"

# Auto-generated cell writes for sheet4 ("options") rows 86-98
$ws4 = $wb.Worksheets.Item("options")
$ws4.Cells.Item(86, 1).Value = "Discap_YTY_opts"
$ws4.Cells.Item(86, 2).Value = "D1100"
$ws4.Cells.Item(86, 3).Value = "(1100 mm rulo dış çapı)"
$ws4.Cells.Item(86, 4).Value = 1
$ws4.Cells.Item(87, 1).Value = "Discap_YTY_opts"
$ws4.Cells.Item(87, 2).Value = "D1200"
$ws4.Cells.Item(87, 3).Value = "(1200 mm rulo dış çapı)"
$ws4.Cells.Item(87, 4).Value = 2
$ws4.Cells.Item(88, 1).Value = "Discap_YTY_opts"
$ws4.Cells.Item(88, 2).Value = "D1300"
$ws4.Cells.Item(88, 3).Value = "(1300 mm rulo dış çapı)"
$ws4.Cells.Item(88, 4).Value = 3
$ws4.Cells.Item(89, 1).Value = "Discap_YTY_opts"
$ws4.Cells.Item(89, 2).Value = "D1400"
$ws4.Cells.Item(89, 3).Value = "(1400 mm rulo dış çapı)"
$ws4.Cells.Item(89, 4).Value = 4
$ws4.Cells.Item(90, 1).Value = "Discap_YTY_opts"
$ws4.Cells.Item(90, 2).Value = "D1500"
$ws4.Cells.Item(90, 3).Value = "(1500 mm rulo dış çapı)"
$ws4.Cells.Item(90, 4).Value = 5
$ws4.Cells.Item(91, 1).Value = "Dogrultma_Tipi_YTY_opts"
$ws4.Cells.Item(91, 2).Value = "T-3092"
$ws4.Cells.Item(91, 3).Value = "Ø30mm - 9 + 2 merdaneli doğrultucu"
$ws4.Cells.Item(91, 4).Value = 1
$ws4.Cells.Item(92, 1).Value = "Hiz_acici_YTY_opts"
$ws4.Cells.Item(92, 2).Value = "VR35"
$ws4.Cells.Item(92, 3).Value = "(35 dev/dk)"
$ws4.Cells.Item(92, 4).Value = 1
$ws4.Cells.Item(93, 1).Value = "Hiz_acici_YTY_opts"
$ws4.Cells.Item(93, 2).Value = "VR50"
$ws4.Cells.Item(93, 3).Value = "(50 dev/dk)"
$ws4.Cells.Item(93, 4).Value = 2
$ws4.Cells.Item(94, 1).Value = "Hiz_acici_YTY_opts"
$ws4.Cells.Item(94, 2).Value = "Diğer"
$ws4.Cells.Item(94, 3).Value = "Lütfen aşağıdaki alana değer giriniz"
$ws4.Cells.Item(94, 4).Value = 3
$ws4.Cells.Item(95, 1).Value = "Hiz_acici_harici_opts"
$ws4.Cells.Item(95, 4).Value = 1
$ws4.Cells.Item(96, 1).Value = "Hiz_dogrultucu_YTY_opts"
$ws4.Cells.Item(96, 2).Value = "VS24"
$ws4.Cells.Item(96, 3).Value = "(24 m/dk hız)"
$ws4.Cells.Item(96, 4).Value = 1
$ws4.Cells.Item(97, 1).Value = "Hiz_dogrultucu_YTY_opts"
$ws4.Cells.Item(97, 2).Value = "Diğer"
$ws4.Cells.Item(97, 3).Value = "Lütfen aşağıdaki alana değer giriniz"
$ws4.Cells.Item(97, 4).Value = 2
$ws4.Cells.Item(98, 1).Value = "Hiz_dogrultucu_harici_opts"
$ws4.Cells.Item(98, 4).Value = 1

# ---------------------------------------------------------------------------
# View state: the author ended the session with the "sections" sheet active
# (selection on D45) after having worked in "fields" (selection left on
# C43) and "options" (selection left on A99, scrolled down to the new rows).
# ---------------------------------------------------------------------------
$ws3.Range("C43").Select()
$ws4.Range("A99").Select()

$ws2 = $wb.Worksheets.Item("sections")
$ws2.Activate()
$ws2.Range("D45").Select()
